$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

# Text columns: force "Text" number format so Excel doesn't auto-convert
# the string to a date serial / number, then clear the formatting so the
# cell ends up with the default (no explicit) style, matching the other
# data rows in this sheet.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-06"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "14:26:11"
$ws.Cells.Item($row, 2).ClearFormats()

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"
$ws.Cells.Item($row, 4).ClearFormats()

# Numeric columns
$ws.Cells.Item($row, 5).Value = 140494
$ws.Cells.Item($row, 6).Value = 143007
$ws.Cells.Item($row, 7).Value = 172058
$ws.Cells.Item($row, 8).Value = 147268
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118141
$ws.Cells.Item($row, 11).Value = 224561
$ws.Cells.Item($row, 12).Value = 248985
$ws.Cells.Item($row, 13).Value = 184936
$ws.Cells.Item($row, 14).Value = 110338
$ws.Cells.Item($row, 15).Value = 40590
$ws.Cells.Item($row, 16).Value = 30830
$ws.Cells.Item($row, 17).Value = 72471
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41877
$ws.Cells.Item($row, 20).Value = -1
